# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.594.46"
$ws.Range("E2").Value = "  +3.60%  "

$ws.Range("D3").Value = "3.203.87"
$ws.Range("E3").Value = "  +1.36%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "638.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.02%  "

$ws.Range("E7").Value = "  +6.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.707"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").ClearFormats()

$ws.Range("D10").Value = "3.202.54"
$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("E11").Value = "  +8.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.181"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.97%  "

$ws.Range("E13").Value = "  +8.33%  "

$ws.Range("E14").Value = "  +4.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.39"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.57%  "

$ws.Range("D16").Value = "90.263.66"
$ws.Range("E16").Value = "  +3.55%  "

$ws.Range("D17").Value = "3.787.47"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").Value = "3.201.67"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000226"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +76.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "437.89"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.34"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.11%  "

$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("E27").Value = "  +11.38%  "

$ws.Range("D28").Value = "3.372.69"
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.993"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +41.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.45"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "541.34"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.08"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.91"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.04%  "

$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.53"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.53%  "

$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("E41").Value = "  -3.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.95"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.41%  "

$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("E44").Value = "  +2.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "146.10"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.77"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "173.32"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("E49").Value = "  +8.94%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.623"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +7.45%  "

$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.87%  "
